$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ D = "27.368.97"; E = "  +1.45%  " }
  3 = @{ D = "1.828.12"; E = "  +0.28%  " }
  4 = @{ E = "  +0.02%  " }
  5 = @{ E = "  +1.09%  " }
  6 = @{ E = "  +0.04%  " }
  7 = @{ D = "0.4476"; E = "  +5.42%  " }
  8 = @{ D = "0.3769"; E = "  +3.13%  " }
  9 = @{ E = "  +4.29%  " }
  10 = @{ D = "0.8949"; E = "  +6.42%  " }
  11 = @{ D = "21.08"; E = "  +2.68%  " }
  12 = @{ D = "1.815.16"; E = "  -0.46%  " }
  13 = @{ D = "6.748"; E = "  +1.60%  " }
  14 = @{ D = "94.62"; E = "  +5.72%  " }
  15 = @{ D = "5.419"; E = "  +2.80%  " }
  16 = @{ D = "0.07125"; E = "  +1.12%  " }
  17 = @{ E = "  -0.05%  " }
  18 = @{ D = "0.000008836"; E = "  +1.06%  " }
  19 = @{ E = "  +0.00%  " }
  20 = @{ D = "15.26"; E = "  +2.90%  " }
  21 = @{ D = "27.389.00"; E = "  +1.23%  " }
  22 = @{ D = "5.290"; E = "  +3.35%  " }
  23 = @{ E = "  +1.63%  " }
  24 = @{ D = "2.005"; E = "  +1.37%  " }
  25 = @{ D = "2.462"; E = "  +11.16%  " }
  26 = @{ D = "151.62"; E = "  +0.58%  " }
  27 = @{ D = "18.67"; E = "  +2.89%  " }
  28 = @{ D = "5.379"; E = "  +3.35%  " }
  29 = @{ D = "118.01"; E = "  +1.28%  " }
  30 = @{ D = "0.08847"; E = "  +1.63%  " }
  31 = @{ D = "0.7863"; E = "  +7.14%  " }
  32 = @{ D = "1.205"; E = "  +2.80%  " }
  33 = @{ D = "4.559"; E = "  +3.42%  " }
  34 = @{ D = "2.892"; E = "  -0.28%  " }
  35 = @{ D = "1.001"; E = "  +0.05%  " }
  36 = @{ D = "1.111"; E = "  +2.11%  " }
  37 = @{ D = "0.01994"; E = "  +3.14%  " }
  38 = @{ D = "0.05339"; E = "  +2.56%  " }
  39 = @{ D = "7.405" }
  40 = @{ D = "0.5333"; E = "  +4.56%  " }
  41 = @{ D = "0.1736"; E = "  +3.17%  " }
  42 = @{ D = "2.873"; E = "  -0.07%  " }
  43 = @{ D = "2.300"; E = "  +17.88%  " }
  44 = @{ D = "8.831"; E = "  +3.82%  " }
  45 = @{ D = "0.5182"; E = "  +9.69%  " }
  46 = @{ D = "10.79"; E = "  +2.18%  " }
  47 = @{ B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "106.50"; E = "  +0.79%  " }
  48 = @{ B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "1.712"; E = "  +3.98%  " }
  49 = @{ D = "1.001"; E = "  +0.13%  " }
  50 = @{ D = "0.06376"; E = "  +0.96%  " }
  51 = @{ D = "64.55"; E = "  +3.50%  " }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ref = "$col$row"
        if ($col -eq "D") {
            $ws.Range($ref).NumberFormat = "@"
            $ws.Range($ref).Value = $vals[$col]
            $ws.Range($ref).ClearFormats()
        } else {
            $ws.Range($ref).Value = $vals[$col]
        }
    }
}